# PROYEK1-LogBookRegawa.xlsx — "Tambah hasil scrap dalam format .csv"
#
# Adds a new log-book entry (row 7 / spreadsheet row 10-11): records the
# "Durasi" + "Nama Kegiatan" for the web-scraping activity, extending the
# existing merged B10:C10 cell down to B10:C11 so the "13 Maret 2018" date
# now spans both new activity rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 (activity #7): duration + activity description -----------------
$ws.Range("D10").Value = "1 jam"
$ws.Range("E10").Value = "Mencari web untuk di scrap dan update list web untuk di scrap"

# --- Extend the date merge B10:C10 down to cover the new row 11 ------------
$ws.Range("B10:C11").Merge()

# --- Row 11 (activity #8): duration + activity description -----------------
$ws.Range("D11").Value = "2 jam"
$ws.Range("E11").Value = "Scraping data dari web www.pumasera.com"

# --- Selection / zoom, matching where the author left off -------------------
[void]$ws.Range("E10").Select()
$excel.ActiveWindow.Zoom = 115
